$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "IAM"

# Delete row 3 (the PUT "Upate user profile" test case) entirely
$ws.Rows.Item(3).Delete()

# Update row 2 cell values to describe the new IAM "authorize" test
$ws.Range("B2").Value = "Test authorize API"
$ws.Range("C2").Value = "1PAUTH"
$ws.Range("D2").Value = "/authorize"
$ws.Range("G2").Value = "?provider=thomsonreuters&backurl=%2Fui%2F%23%2Flogin%2FAUTHTOKEN"
$ws.Range("J2").Value = "status=200"
$ws.Range("K2").Value = ""

# Give G2 the same border+fill-applied style used elsewhere (applyFill + applyBorder, no visible border/fill)
$ws.Range("G2").Interior.Color = $null
$ws.Range("G2").Borders.Color = $null

# Re-fit column widths to the new (longer/shorter) content
$ws.Columns.Item(2).ColumnWidth = 27.1666666666667
$ws.Columns.Item(4).ColumnWidth = 79.42
$ws.Columns.Item(7).ColumnWidth = 69.59

# Update the sheet view / selection to match
$ws.Range("A2").Select()
